$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Dati Input"
$ws2 = $wb.Worksheets.Item(2)   # "Dati Output"

# --- Sheet 1: "Dati Input" ------------------------------------------------
# Copy column J's formatting into the new column K for every used row so the
# new cells inherit the same row styles (header / ROOT / body) as the rest
# of the table, then fill in the "Caso Test" values.
$ws1.Range("J2:J22").Copy()
$ws1.Range("K2:K22").PasteSpecial(-4122)

$ws1.Range("K2").Value = "Caso Test"

$ws1.Range("K4").Value = "abc"
$ws1.Range("K5").Value = "abc"

$ws1.Range("K7").Value = "efd"
$ws1.Range("K8").Value = "efd"

$ws1.Range("K10").Value = "efd"
$ws1.Range("K11").Value = "efd"

$ws1.Range("K13").Value = 123
$ws1.Range("K14").Value = 123

$ws1.Range("K16").Value = "ghi"
$ws1.Range("K17").Value = 456

$ws1.Range("K19").Value = "jkl"
$ws1.Range("K21").Value = "jkl"
$ws1.Range("K22").Value = "jkl"

# --- Sheet 2: "Dati Output" ------------------------------------------------
$ws2.Range("J2:J11").Copy()
$ws2.Range("K2:K11").PasteSpecial(-4122)

$ws2.Range("K2").Value = "Caso Test"

$ws2.Range("K7").Value = "abc"
$ws2.Range("K8").Value = "abc"
$ws2.Range("K9").Value = "abc"

$ws2.Range("K10").Value = "efd"

$ws2.Range("K11").Value = 123

# --- Selections / active sheet --------------------------------------------
# Final state: "Dati Input" selection is K2:K22 (not the active tab);
# "Dati Output" is the active tab with L5 selected.
$ws1.Activate()
$ws1.Range("K2:K22").Select()

$ws2.Activate()
$ws2.Range("L5").Select()
